$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the column headers:
#      *_old -> *_FV2304   (columns A..J)
#      diff stays the same (column K)
#      *_new -> *_FV2310   (columns L..U)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2304"
$ws.Range("B1").Value = "Segmentgruppe_FV2304"
$ws.Range("C1").Value = "Segment_FV2304"
$ws.Range("D1").Value = "Datenelement_FV2304"
$ws.Range("E1").Value = "Segment ID_FV2304"
$ws.Range("F1").Value = "Code_FV2304"
$ws.Range("G1").Value = "Qualifier_FV2304"
$ws.Range("H1").Value = "Beschreibung_FV2304"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("J1").Value = "Bedingung_FV2304"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2310"
$ws.Range("M1").Value = "Segmentgruppe_FV2310"
$ws.Range("N1").Value = "Segment_FV2310"
$ws.Range("O1").Value = "Datenelement_FV2310"
$ws.Range("P1").Value = "Segment ID_FV2310"
$ws.Range("Q1").Value = "Code_FV2310"
$ws.Range("R1").Value = "Qualifier_FV2310"
$ws.Range("S1").Value = "Beschreibung_FV2310"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("U1").Value = "Bedingung_FV2310"

# ---------------------------------------------------------------------------
# 2) Turn the data range into an Excel Table ("Table1") spanning A1:U65.
#    We first build the table on an empty, unformatted scratch area so the
#    already-bold/filled header row (row 1) does not get baked into a new
#    dxf / headerRowDxfId when the table adopts it. After that we resize the
#    table onto the real data range and clean the scratch cells again.
# ---------------------------------------------------------------------------
$scratchHeader = $ws.Range("A67:B67")
$scratchData = $ws.Range("A68:B68")
$scratchHeader.Value = @("tmp1", "tmp2")
$scratchData.Value = @("tmp_a", "tmp_b")

$tbl = $ws.ListObjects.Add(1, $ws.Range("A67:B68"), $null, 1)
$tbl.Name = "Table1"

$ws.Range("A67:B68").Clear()

$tbl.Resize($ws.Range("A1:U65"))

# Re-apply the header names once more, now bound to the table, so the
# table columns pick up the final (already-written) header text.
$ws.Range("A1").Value = "Segmentname_FV2304"
$ws.Range("B1").Value = "Segmentgruppe_FV2304"
$ws.Range("C1").Value = "Segment_FV2304"
$ws.Range("D1").Value = "Datenelement_FV2304"
$ws.Range("E1").Value = "Segment ID_FV2304"
$ws.Range("F1").Value = "Code_FV2304"
$ws.Range("G1").Value = "Qualifier_FV2304"
$ws.Range("H1").Value = "Beschreibung_FV2304"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("J1").Value = "Bedingung_FV2304"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2310"
$ws.Range("M1").Value = "Segmentgruppe_FV2310"
$ws.Range("N1").Value = "Segment_FV2310"
$ws.Range("O1").Value = "Datenelement_FV2310"
$ws.Range("P1").Value = "Segment ID_FV2310"
$ws.Range("Q1").Value = "Code_FV2310"
$ws.Range("R1").Value = "Qualifier_FV2310"
$ws.Range("S1").Value = "Beschreibung_FV2310"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("U1").Value = "Bedingung_FV2310"

# Drop the automatically assigned table style name so no extra style
# formatting (dxf) is introduced.
$tbl.TableStyle = ""

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1) and select the pane.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("A1").Select()
